$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "31.03.2026"
$ws.Range("B4").Value = "10:00"
$ws.Range("C4").Value = "55NM123"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = "11:30"
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = "13:45"
$ws.Range("H4").Value = 14
$ws.Range("I4").Value = "14:30"
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = "Nisa Karaman"
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = 10
